$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '34.387.64'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +12.53%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.826.45'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +9.33%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.13%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '230.07'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +5.01%  '
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  -0.15%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '31.52'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +8.44%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '46.83'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +4.25%  '
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +9.61%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0678'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +6.40%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.0932'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +3.19%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '2.090.68'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +9.40%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '1.829.37'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +9.42%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.653'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +8.34%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '34.370.91'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +12.54%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '10.28'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +3.55%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '4.31'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +7.85%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '70.38'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '257.72'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +6.61%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.0₃0758'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +5.50%  '
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -0.24%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '10.66'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +7.33%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '4.34'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +3.04%  '
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +3.87%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '159.34'
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +0.34%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '16.77'
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  +6.53%  '
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +5.35%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '7.17'
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +7.95%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -0.14%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '3.89'
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +13.00%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.0524'
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +6.61%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.21'
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +6.37%  '
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +9.24%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.538.60'
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +2.92%  '
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +6.09%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.638'
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +7.56%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '84.37'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +1.64%  '
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +5.39%  '
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +2.43%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.915'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +9.68%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.12'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +6.10%  '
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +6.05%  '
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +6.19%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.980.60'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +9.56%  '
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +5.84%  '
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +18.49%  '
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -0.26%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '51.63'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  +4.73%  '
